$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "bleu" = "noir"
    "pas de résultat ni de publication" = "pas de résultat postés ni publiés"
    "résultat et / ou publication posté dans les 36 mois" = "résultat postés ou publiés dans les 36 mois"
    "résultat et / ou publication posté" = "résultat postés ou publiés"
}

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($map.ContainsKey($val)) {
            $cell.Value = $map[$val]
        }
    }
}
